$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63, shifting existing rows 63-213 down to 64-214
$ws.Rows(63).Insert()

# Populate the newly inserted row 63 with the new data record
$ws.Cells.Item(63, 1).Value = 9
$ws.Cells.Item(63, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(63, 3).Value = "Metropolitana"
$ws.Cells.Item(63, 4).Value = 44581
$ws.Cells.Item(63, 5).Value = 13
$ws.Cells.Item(63, 6).Value = 300000001
$ws.Cells.Item(63, 7).Value = "Rabanito"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 6100
$ws.Cells.Item(63, 11).Value = 3000
$ws.Cells.Item(63, 12).Value = 3000
$ws.Cells.Item(63, 13).Value = 3000
$ws.Cells.Item(63, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(63, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(63, 16).Value = 30
$ws.Cells.Item(63, 17).Value = 100
$ws.Cells.Item(63, 18).Value = "Hortaliza"
